$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.6122626666666666
$ws.Range("H2").Value = 1.836788
$ws.Range("I2").Value = 0.006779070576782467
$ws.Range("J2").Value = 0.006779070576782467
$ws.Range("M2").Value = 8.306580666666667
$ws.Range("N2").Value = 24.919742
$ws.Range("O2").Value = 0.4741050717515609
$ws.Range("P2").Value = 0.4741050717515609
$ws.Range("Q2").Value = 5.085809229855111
$ws.Range("R2").Value = 45.772283068696
$ws.Range("S2").Value = 0.003213991742214347
$ws.Range("T2").Value = 0.003213991742214347
$ws.Range("G3").Value = 0.6122626666666666
$ws.Range("H3").Value = 1.836788
$ws.Range("I3").Value = 0.006779070576782467
$ws.Range("J3").Value = 0.006779070576782467
$ws.Range("O3").Value = 0.2729564927611473
$ws.Range("P3").Value = 0.2729564927611473
$ws.Range("Q3").Value = 2.928052731232889
$ws.Range("R3").Value = 26.35247458109599
$ws.Range("S3").Value = 0.00185039132881883
$ws.Range("T3").Value = 0.00185039132881883
$ws.Range("G4").Value = 0.6122626666666666
$ws.Range("H4").Value = 1.836788
$ws.Range("I4").Value = 0.006779070576782467
$ws.Range("J4").Value = 0.006779070576782467
$ws.Range("M4").Value = 3.901832333333334
$ws.Range("N4").Value = 11.705497
$ws.Range("O4").Value = 0.2227003592201187
$ws.Range("P4").Value = 0.2227003592201188
$ws.Range("Q4").Value = 2.388946269292889
$ws.Range("R4").Value = 21.500516423636
$ws.Range("S4").Value = 0.001509701452627993
$ws.Range("T4").Value = 0.001509701452627993
$ws.Range("G5").Value = 0.6122626666666666
$ws.Range("H5").Value = 1.836788
$ws.Range("I5").Value = 0.006779070576782467
$ws.Range("J5").Value = 0.006779070576782467
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5297876666666667
$ws.Range("N5").Value = 1.589363
$ws.Range("O5").Value = 0.03023807626717307
$ws.Range("P5").Value = 0.03023807626717307
$ws.Range("Q5").Value = 0.3243692095604444
$ws.Range("R5").Value = 2.919322886044
$ws.Range("S5").Value = 0.0002049860531212972
$ws.Range("T5").Value = 0.0002049860531212972
$ws.Range("I6").Value = 0.003538518590750013
$ws.Range("J6").Value = 0.003538518590750013
$ws.Range("M6").Value = 8.306580666666667
$ws.Range("N6").Value = 24.919742
$ws.Range("O6").Value = 0.4741050717515609
$ws.Range("P6").Value = 0.4741050717515609
$ws.Range("Q6").Value = 2.654675195518
$ws.Range("R6").Value = 23.892076759662
$ws.Range("S6").Value = 0.001677629610361767
$ws.Range("T6").Value = 0.001677629610361767
$ws.Range("I7").Value = 0.003538518590750013
$ws.Range("J7").Value = 0.003538518590750013
$ws.Range("O7").Value = 0.2729564927611473
$ws.Range("P7").Value = 0.2729564927611473
$ws.Range("S7").Value = 0.000965861624101241
$ws.Range("T7").Value = 0.000965861624101241
$ws.Range("I8").Value = 0.003538518590750013
$ws.Range("J8").Value = 0.003538518590750013
$ws.Range("M8").Value = 3.901832333333334
$ws.Range("N8").Value = 11.705497
$ws.Range("O8").Value = 0.2227003592201187
$ws.Range("P8").Value = 0.2227003592201188
$ws.Range("Q8").Value = 1.246974889913
$ws.Range("R8").Value = 11.222774009217
$ws.Range("S8").Value = 0.0007880293612670963
$ws.Range("T8").Value = 0.0007880293612670964
$ws.Range("I9").Value = 0.003538518590750013
$ws.Range("J9").Value = 0.003538518590750013
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5297876666666667
$ws.Range("N9").Value = 1.589363
$ws.Range("O9").Value = 0.03023807626717307
$ws.Range("P9").Value = 0.03023807626717307
$ws.Range("Q9").Value = 0.169313251027
$ws.Range("R9").Value = 1.523819259243
$ws.Range("S9").Value = 0.0001069979950199087
$ws.Range("T9").Value = 0.0001069979950199087
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.115957
$ws.Range("H10").Value = 0.347871
$ws.Range("I10").Value = 0.001283894527085267
$ws.Range("J10").Value = 0.001283894527085267
$ws.Range("M10").Value = 8.306580666666667
$ws.Range("N10").Value = 24.919742
$ws.Range("O10").Value = 0.4741050717515609
$ws.Range("P10").Value = 0.4741050717515609
$ws.Range("Q10").Value = 0.9632061743646666
$ws.Range("R10").Value = 8.668855569282
$ws.Range("S10").Value = 0.0006087009068851968
$ws.Range("T10").Value = 0.0006087009068851969
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.115957
$ws.Range("H11").Value = 0.347871
$ws.Range("I11").Value = 0.001283894527085267
$ws.Range("J11").Value = 0.001283894527085267
$ws.Range("O11").Value = 0.2729564927611473
$ws.Range("P11").Value = 0.2729564927611473
$ws.Range("Q11").Value = 0.5545466497313333
$ws.Range("R11").Value = 4.990919847581999
$ws.Range("S11").Value = 0.0003504473471884263
$ws.Range("T11").Value = 0.0003504473471884263
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.115957
$ws.Range("H12").Value = 0.347871
$ws.Range("I12").Value = 0.001283894527085267
$ws.Range("J12").Value = 0.001283894527085267
$ws.Range("M12").Value = 3.901832333333334
$ws.Range("N12").Value = 11.705497
$ws.Range("O12").Value = 0.2227003592201187
$ws.Range("P12").Value = 0.2227003592201188
$ws.Range("Q12").Value = 0.4524447718763333
$ws.Range("R12").Value = 4.072002946887
$ws.Range("S12").Value = 0.0002859237723826334
$ws.Range("T12").Value = 0.0002859237723826335
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.115957
$ws.Range("H13").Value = 0.347871
$ws.Range("I13").Value = 0.001283894527085267
$ws.Range("J13").Value = 0.001283894527085267
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.5297876666666667
$ws.Range("N13").Value = 1.589363
$ws.Range("O13").Value = 0.03023807626717307
$ws.Range("P13").Value = 0.03023807626717307
$ws.Range("Q13").Value = 0.06143258846366666
$ws.Range("R13").Value = 0.552893296173
$ws.Range("S13").Value = 0.0000388225006290104
$ws.Range("T13").Value = 0.00003882250062901042
$ws.Range("G14").Value = 89.26880233333334
$ws.Range("H14").Value = 267.806407
$ws.Range("I14").Value = 0.9883985163053822
$ws.Range("J14").Value = 0.9883985163053823
$ws.Range("M14").Value = 8.306580666666667
$ws.Range("N14").Value = 24.919742
$ws.Range("O14").Value = 0.4741050717515609
$ws.Range("P14").Value = 0.4741050717515609
$ws.Range("Q14").Value = 741.5185075985549
$ws.Range("R14").Value = 6673.666568386995
$ws.Range("S14").Value = 0.4686047494920995
$ws.Range("T14").Value = 0.4686047494920996
$ws.Range("G15").Value = 89.26880233333334
$ws.Range("H15").Value = 267.806407
$ws.Range("I15").Value = 0.9883985163053822
$ws.Range("J15").Value = 0.9883985163053823
$ws.Range("O15").Value = 0.2729564927611473
$ws.Range("P15").Value = 0.2729564927611473
$ws.Range("Q15").Value = 426.9144187886771
$ws.Range("R15").Value = 3842.229769098094
$ws.Range("S15").Value = 0.2697897924610387
$ws.Range("T15").Value = 0.2697897924610388
$ws.Range("G16").Value = 89.26880233333334
$ws.Range("H16").Value = 267.806407
$ws.Range("I16").Value = 0.9883985163053822
$ws.Range("J16").Value = 0.9883985163053823
$ws.Range("M16").Value = 3.901832333333334
$ws.Range("N16").Value = 11.705497
$ws.Range("O16").Value = 0.2227003592201187
$ws.Range("P16").Value = 0.2227003592201188
$ws.Range("Q16").Value = 348.3118993021421
$ws.Range("R16").Value = 3134.80709371928
$ws.Range("S16").Value = 0.220116704633841
$ws.Range("T16").Value = 0.2201167046338411
$ws.Range("G17").Value = 89.26880233333334
$ws.Range("H17").Value = 267.806407
$ws.Range("I17").Value = 0.9883985163053822
$ws.Range("J17").Value = 0.9883985163053823
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.5297876666666667
$ws.Range("N17").Value = 1.589363
$ws.Range("O17").Value = 0.03023807626717307
$ws.Range("P17").Value = 0.03023807626717307
$ws.Range("Q17").Value = 47.29351049430456
$ws.Range("R17").Value = 425.6415944487411
$ws.Range("S17").Value = 0.02988726971840286
$ws.Range("T17").Value = 0.02988726971840286
